# SwaadSutra order sheet: a new order came in (Sagar Borse, Jawar Bhakari x1)
# at 2026-01-13 18:59. New orders are prepended right under the header, so
# insert a fresh row at row 2 and push all the existing orders down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Orders")
$summary = $wb.Worksheets.Item("Daily Summary")

$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).Value = 8                       # Order ID
$ws.Cells.Item(2, 2).Value = "2026-01-13 18:59"       # Date
$ws.Cells.Item(2, 3).Value = "Sagar Borse"             # Customer
$ws.Cells.Item(2, 4).Value = "A-1608"                  # Flat No

# Phone must stay text (not be coerced to a number) - force text format first.
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "7588930329"              # Phone

$ws.Cells.Item(2, 6).Value = "Jawar Bhakari x1"        # Items
$ws.Cells.Item(2, 7).Value = 20                        # Total
$ws.Cells.Item(2, 8).Value = "NEW"                     # Status
$ws.Cells.Item(2, 9).Value = "PENDING"                 # Payment

# Collection Date must stay text (not be coerced to a date) too.
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "2026-01-16"             # Collection Date

$ws.Cells.Item(2, 11).Value = "10:00"                  # Collection Time
# Notes / Cancel Reason / Feedback (L2:N2) are left blank for the new order.

# Daily Summary for 2026-01-13: one more order, +20 to Revenue and Pending.
$summary.Cells.Item(2, 2).Value = 8     # Total Orders
$summary.Cells.Item(2, 5).Value = 215   # Revenue
$summary.Cells.Item(2, 7).Value = 215   # Pending
